$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Current Asymptomatic Infections" (D) and "Current Infectious Carriers" (E)
# for rows 9 through 13 to reflect pushed-back start date / updated infection stats.
$ws.Range("D9").Value = 26
$ws.Range("E9").Value = 8

$ws.Range("D10").Value = 32
$ws.Range("E10").Value = 14

$ws.Range("D11").Value = 38
$ws.Range("E11").Value = 20

$ws.Range("D12").Value = 44
$ws.Range("E12").Value = 26

$ws.Range("D13").Value = 48
$ws.Range("E13").Value = 30
